# Update simulation results in the Weights workbook after re-running the
# aircraft model tests ("More tests on the new aircraft models.")
#
# Each worksheet below gets a handful of numeric cells refreshed with the
# latest computed values (mass/weight results in column C, and percentage
# deltas in column D where present).

$wb = $excel.ActiveWorkbook

# ---- GLOBAL RESULTS --------------------------------------------------
$ws = $wb.Worksheets.Item("GLOBAL RESULTS")
$ws.Range("C2").Value  = 1100.0
$ws.Range("C6").Value  = 33735.068110578824
$ws.Range("C7").Value  = 33735.068110578824
$ws.Range("C8").Value  = 32723.01606726146
$ws.Range("C11").Value = 5784.462396351928
$ws.Range("C13").Value = 27950.605714226906
$ws.Range("C14").Value = 27950.605714226906
$ws.Range("C15").Value = 18410.605714226906
$ws.Range("C16").Value = 18104.5475202269
$ws.Range("C17").Value = 17329.01752022691
$ws.Range("C19").Value = 0.0
$ws.Range("C21").Value = 330828.0056866077
$ws.Range("C22").Value = 330828.0056866077
$ws.Range("C23").Value = 320903.16551600944
$ws.Range("C27").Value = 274101.80752742314
$ws.Range("C28").Value = 274101.80752742314
$ws.Range("C29").Value = 180546.3665274232
$ws.Range("C30").Value = 177544.96093923307
$ws.Range("C31").Value = 169939.6096647332
$ws.Range("C33").Value = 0.0

# ---- FUSELAGE ----------------------------------------------------------
$ws = $wb.Worksheets.Item("FUSELAGE")
$ws.Range("C7").Value  = 3909.0
$ws.Range("D7").Value  = 1.00775193798452
$ws.Range("C8").Value  = 3833.0
$ws.Range("D8").Value  = -0.9560723514211653
$ws.Range("C9").Value  = 3356.0
$ws.Range("D9").Value  = -13.281653746770006
$ws.Range("C12").Value = 4512.666666666666
$ws.Range("D12").Value = 16.606373815676168

# ---- WING ----------------------------------------------------------------
$ws = $wb.Worksheets.Item("WING")
$ws.Range("C7").Value  = 3021.0
$ws.Range("D7").Value  = 17.093023255813993
$ws.Range("C13").Value = 2630.7142857142853
$ws.Range("D13").Value = 1.965669988925823

# ---- HORIZONTAL TAIL -------------------------------------------------
$ws = $wb.Worksheets.Item("HORIZONTAL TAIL")
$ws.Range("C8").Value  = 177.0
$ws.Range("D8").Value  = -54.26356589147286
$ws.Range("C9").Value  = 185.0
$ws.Range("D9").Value  = -52.19638242894055
$ws.Range("C10").Value = 226.66666666666663
$ws.Range("D10").Value = -41.42980189491815

# ---- POWER PLANT -------------------------------------------------------
$ws = $wb.Worksheets.Item("POWER PLANT")
$ws.Range("C2").Value  = 3507.999999999999
$ws.Range("C3").Value  = 4830.515999999998
$ws.Range("C8").Value  = 1754.0
$ws.Range("C9").Value  = 2415.257999999999
$ws.Range("C12").Value = 1754.0
$ws.Range("C13").Value = 2415.257999999999

# ---- LANDING GEARS -----------------------------------------------------
$ws = $wb.Worksheets.Item("LANDING GEARS")
$ws.Range("C5").Value = 1053.0
$ws.Range("D5").Value = 2.0348837209302775
$ws.Range("C6").Value = 1349.0
$ws.Range("D6").Value = 30.71705426356595
$ws.Range("C7").Value = 1521.0
$ws.Range("D7").Value = 47.38372093023262
$ws.Range("C8").Value = 1367.0
$ws.Range("D8").Value = 32.46124031007758
$ws.Range("C9").Value = 1322.5
$ws.Range("D9").Value = 28.149224806201573
